$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.229822
$ws.Range("N2").Value = 0.689466
$ws.Range("O2").Value = 0.09226175421862418
$ws.Range("P2").Value = 0.09226175421862419
$ws.Range("Q2").Value = 0.004538295034
$ws.Range("R2").Value = 0.040844655306
$ws.Range("S2").Value = 0.06589124269638957
$ws.Range("T2").Value = 0.06589124269638957

# Row 3
$ws.Range("O3").Value = 0.4364142651333466
$ws.Range("P3").Value = 0.4364142651333466
$ws.Range("S3").Value = 0.3116771245420685
$ws.Range("T3").Value = 0.3116771245420685

# Row 4
$ws.Range("M4").Value = 1.174057666666666
$ws.Range("O4").Value = 0.4713239806480292
$ws.Range("P4").Value = 0.4713239806480293
$ws.Range("S4").Value = 0.3366088479514153
$ws.Range("T4").Value = 0.3366088479514154

# Row 5
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.229822
$ws.Range("N5").Value = 0.689466
$ws.Range("O5").Value = 0.09226175421862418
$ws.Range("P5").Value = 0.09226175421862419
$ws.Range("Q5").Value = 0.001816283266
$ws.Range("R5").Value = 0.016346549394
$ws.Range("S5").Value = 0.02637051152223461
$ws.Range("T5").Value = 0.02637051152223462

# Row 6
$ws.Range("O6").Value = 0.4364142651333466
$ws.Range("P6").Value = 0.4364142651333466
$ws.Range("Q6").Value = 0.008591338128333334
$ws.Range("R6").Value = 0.07732204315499999
$ws.Range("S6").Value = 0.1247371405912781
$ws.Range("T6").Value = 0.1247371405912781

# Row 7
$ws.Range("M7").Value = 1.174057666666666
$ws.Range("O7").Value = 0.4713239806480292
$ws.Range("P7").Value = 0.4713239806480293
$ws.Range("Q7").Value = 0.009278577739666666
$ws.Range("R7").Value = 0.08350719965699999
$ws.Range("S7").Value = 0.1347151326966139
$ws.Range("T7").Value = 0.134715132696614
